$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mouser")
$ws.Rows("4:13").Delete()
